$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 11
$ws.Range("C12").Value = 20

$ws.Range("A13").Select()
